# Generate Report for Archive
# - Update the localization status text from "Ready for handoff" to
#   "In Translation" everywhere it appears (Overview!E2:F4, and the
#   "Status" column on each per-language sheet).
# - The status columns re-size (narrower, since "In Translation" is
#   shorter than "Ready for handoff") to fit the new content.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# Overview sheet: zh-cn/de-de status columns (E and F) for each file row.
$wsOverview.Range("E2:F4").Value = $newStatus

# Per-language detail sheets: "Status" column (C) for each file row.
$wsZhCn.Range("C2:C4").Value = $newStatus
$wsDeDe.Range("C2:C4").Value = $newStatus

# The status columns auto-size to fit the (now shorter) content.
$newColumnWidth = 12.5
$wsOverview.Columns("E:F").ColumnWidth = $newColumnWidth
$wsZhCn.Columns("C:C").ColumnWidth = $newColumnWidth
$wsDeDe.Columns("C:C").ColumnWidth = $newColumnWidth
